# Auto-generated edit script: updates crypto price/volume table cells
# per the commit diff (row reorders + value refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.825.01'
$ws.Range('E2').Value = '  +0.31%  '

$ws.Range('D3').Value = '2.534.95'
$ws.Range('E3').Value = '  +0.82%  '

$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.64'
$ws.Range('E5').Value = '  -0.09%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '95.81'
$ws.Range('E6').Value = '  -0.01%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.575'

$ws.Range('E8').Value = '  +0.06%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.533'
$ws.Range('E9').Value = '  -1.11%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.95'
$ws.Range('E10').Value = '  -0.31%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0808'
$ws.Range('E11').Value = '  -0.61%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.56'
$ws.Range('E12').Value = '  -0.49%  '

$ws.Range('E13').Value = '  -2.57%  '

$ws.Range('D14').Value = '2.928.20'
$ws.Range('E14').Value = '  +0.98%  '

$ws.Range('D15').Value = '2.554.95'
$ws.Range('E15').Value = '  +1.32%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.16'
$ws.Range('E16').Value = '  -2.55%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.850'
$ws.Range('E17').Value = '  -1.31%  '

$ws.Range('D18').Value = '42.891.54'
$ws.Range('E18').Value = '  +0.51%  '

$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.82'
$ws.Range('E19').Value = '  +3.48%  '

$ws.Range('B20').Value = 'InternetComputer(DFINITY)'
$ws.Range('C20').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.93'
$ws.Range('E20').Value = '  +0.13%  '

$ws.Range('D21').Value = '0.0₃0965'
$ws.Range('E21').Value = '  -0.93%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '69.96'
$ws.Range('E22').Value = '  -2.18%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '252.40'
$ws.Range('E23').Value = '  -0.31%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.93'
$ws.Range('E24').Value = '  -1.76%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.07'
$ws.Range('E25').Value = '  +1.33%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.53'
$ws.Range('E26').Value = '  -1.88%  '

$ws.Range('E27').Value = '  +0.04%  '

$ws.Range('E28').Value = '  +1.93%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '40.54'
$ws.Range('E29').Value = '  +6.09%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.41'
$ws.Range('E30').Value = '  +2.62%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.90'
$ws.Range('E31').Value = '  -0.35%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '157.16'
$ws.Range('E32').Value = '  +1.06%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.16'
$ws.Range('E33').Value = '  +3.14%  '

$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.70'
$ws.Range('E34').Value = '  +3.46%  '

$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.32'
$ws.Range('E35').Value = '  -0.87%  '

$ws.Range('B36').Value = 'Celestia'
$ws.Range('C36').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.91'
$ws.Range('E36').Value = '  -4.50%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0780'
$ws.Range('E37').Value = '  -0.96%  '

$ws.Range('E38').Value = '  -1.06%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.118'
$ws.Range('E39').Value = '  -1.43%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.33'
$ws.Range('E40').Value = '  +13.78%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '22.64'
$ws.Range('E41').Value = '  -7.76%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.83'
$ws.Range('E42').Value = '  -0.67%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0303'
$ws.Range('E43').Value = '  -0.06%  '

$ws.Range('E44').Value = '  +0.24%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.29'
$ws.Range('E45').Value = '  -2.78%  '

$ws.Range('D46').Value = '2.028.72'
$ws.Range('E46').Value = '  -0.24%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.12'
$ws.Range('E47').Value = '  +1.80%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '84.71'
$ws.Range('E48').Value = '  +0.14%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '106.00'
$ws.Range('E49').Value = '  +4.30%  '

$ws.Range('D50').Value = '2.780.38'
$ws.Range('E50').Value = '  +0.86%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '74.63'
$ws.Range('E51').Value = '  +0.72%  '
